$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert a new record as a new row 12, pushing the existing
# rows 12-27 down to 13-28 (dimension grows from A1:R27 to A1:R28).
$ws.Rows.Item(12).Insert()

# Fill the newly inserted row 12 with the new weekly price record.
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(12, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(12, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(12, 4).Value = 44540
$ws.Cells.Item(12, 5).Value = 15
$ws.Cells.Item(12, 6).Value = 100112052
$ws.Cells.Item(12, 7).Value = 'Albahaca'
$ws.Cells.Item(12, 8).Value = 'Sin especificar'
$ws.Cells.Item(12, 9).Value = 'Primera'
$ws.Cells.Item(12, 10).Value = 200
$ws.Cells.Item(12, 11).Value = 900
$ws.Cells.Item(12, 12).Value = 1000
$ws.Cells.Item(12, 13).Value = 950
$ws.Cells.Item(12, 14).Value = '$/paquete'
$ws.Cells.Item(12, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(12, 16).Value = 950
$ws.Cells.Item(12, 17).Value = 1
$ws.Cells.Item(12, 18).Value = 'Hortaliza'
